$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 42622.891469907408
$ws.Range("A15").NumberFormat = "m/d/yy h:mm"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 50
$ws.Range("D15").Value = 46
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 99
$ws.Range("G15").Value = 21189
$ws.Range("H15").Value = 11310
$ws.Range("I15").Value = 577
$ws.Range("J15").Value = 137
$ws.Range("K15").Value = 127
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 3
$ws.Range("N15").Value = "Named"
